$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Librero")

# --- Row 132: Cuentos completos / Vladimir Nabokov -----------------------
# Column C written before column B so the shared-string table gets
# "Vladimir Nabokov" (335) before "Cuentos completos" (336), matching the
# author's original entry order.
$ws.Cells.Item(132, 3).Value = "Vladimir Nabokov"
$ws.Cells.Item(132, 2).Value = "Cuentos completos"
$ws.Cells.Item(132, 4).Value = "De bolsillo"
$ws.Cells.Item(132, 5).Value = "Ficción/Contemporánea"
$ws.Cells.Item(132, 7).Value = 0
$ws.Cells.Item(132, 8).Value = 1
$ws.Cells.Item(132, 10).Value = 0
$ws.Cells.Item(132, 11).Value = "Español"
$ws.Cells.Item(132, 12).Value = 1
$ws.Cells.Item(132, 14).Value = 4
$ws.Cells.Item(132, 15).Value = 2021

# --- Row 133: Amistad de Juventud / Alice Munro ---------------------------
$ws.Cells.Item(133, 2).Value = "Amistad de Juventud"
$ws.Cells.Item(133, 3).Value = "Alice Munro"
$ws.Cells.Item(133, 4).Value = "De bolsillo"
$ws.Cells.Item(133, 5).Value = "Ficción/Contemporánea"
$ws.Cells.Item(133, 7).Value = 0
$ws.Cells.Item(133, 8).Value = 1
$ws.Cells.Item(133, 10).Value = 0
$ws.Cells.Item(133, 11).Value = "Español"
$ws.Cells.Item(133, 12).Value = 0
$ws.Cells.Item(133, 15).Value = 2021

# --- Row 134: Los restos del día / Kazuo Ishiguro / Vintage ---------------
$ws.Cells.Item(134, 2).Value = "Los restos del día"
$ws.Cells.Item(134, 4).Value = "Vintage"
$ws.Cells.Item(134, 3).Value = "Kazuo Ishiguro"
$ws.Cells.Item(134, 5).Value = "Ficción/Contemporánea"
$ws.Cells.Item(134, 7).Value = 0
$ws.Cells.Item(134, 8).Value = 0
$ws.Cells.Item(134, 10).Value = 0
$ws.Cells.Item(134, 11).Value = "Español"
$ws.Cells.Item(134, 12).Value = 1
$ws.Cells.Item(134, 14).Value = 5
$ws.Cells.Item(134, 15).Value = 2021

# --- Row 135: El ruiseñor --------------------------------------------------
$ws.Cells.Item(135, 2).Value = "El ruiseñor"
$ws.Cells.Item(135, 4).Value = "De bolsillo"
$ws.Cells.Item(135, 5).Value = "Ficción histórica"
$ws.Cells.Item(135, 7).Value = 0
$ws.Cells.Item(135, 8).Value = 0
$ws.Cells.Item(135, 10).Value = 0
$ws.Cells.Item(135, 11).Value = "Español"
$ws.Cells.Item(135, 12).Value = 0
$ws.Cells.Item(135, 15).Value = 2021

# --- Row 136: Heartstopper / Alice Osimore ---------------------------------
$ws.Cells.Item(136, 2).Value = "Heartstopper"
$ws.Cells.Item(136, 3).Value = "Alice Osimore"

# --- Row 137: chart caption, wrapped across two lines -----------------------
$ws.Cells.Item(137, 2).Value = "Yanely Luna `nGutiérrez"
$ws.Cells.Item(137, 2).WrapText = $true
$ws.Rows.Item(137).RowHeight = 30

# --- Selection bookkeeping on the "Librero" sheet --------------------------
$ws.Range("Q18").Select()

# --- Switch the active tab to "Audiolibros" and set its selection ---------
$ws2 = $wb.Worksheets.Item("Audiolibros")
$ws2.Activate()
$ws2.Range("C9").Select()
